# Kaman new UI - header & footer: add new verification steps to the
# TC08_VerifyProductListPLP sheet (two new rows for a "MountedBearings"
# element, split into EleType1 / EleType2 lookups) and two supporting
# rows on the Testdata sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC08_VerifyProductListPLP")
$ws2 = $wb.Worksheets.Item("Testdata")

# ---------------------------------------------------------------
# Sheet1: insert two new rows before the existing "Ball" row (row 6),
# pushing the old rows 6,7,8 down to 8,9,10.
# ---------------------------------------------------------------
$ws1.Rows.Item(6).Insert()
$ws1.Rows.Item(6).Insert()

# Copy the formatting of the row above (row 5) onto the two new rows
# so the new cells keep the same bordered look as the rest of the table.
$ws1.Range("A5:E5").Copy($ws1.Range("A6:E6"))
$ws1.Range("A5:E5").Copy($ws1.Range("A7:E7"))

# New row 6: WAIT step (only the Keyword column is populated)
$ws1.Range("B6").Value = "WAIT"
$ws1.Range("C6").Value = ""
$ws1.Range("D6").Value = ""
$ws1.Range("E6").Value = ""

# New row 7: MOUSEOVER MountedBearings (CSS)
$ws1.Range("B7").Value = "MOUSEOVER"
$ws1.Range("C7").Value = "MountedBearings"
$ws1.Range("D7").Value = "CSS"
$ws1.Range("E7").Value = ""

# ---------------------------------------------------------------
# Sheet2 (Testdata): append EleType1 / EleType2 -> JSElement rows
# ---------------------------------------------------------------
# Borrow the bordered look already used on B2 for the 4 new cells.
$ws2.Range("B2").Copy($ws2.Range("A4"))
$ws2.Range("B2").Copy($ws2.Range("B4"))
$ws2.Range("B2").Copy($ws2.Range("A5"))
$ws2.Range("B2").Copy($ws2.Range("B5"))

$ws2.Range("A4").Value = "EleType1"
$ws2.Range("B4").Value = "JSElement"
$ws2.Range("A5").Value = "EleType2"
$ws2.Range("B5").Value = "JSElement"

$ws2.Range("A4:B5").Select()

# Re-activate sheet1 and restore its selection last, so it stays the
# active/selected tab when the workbook is saved (matches the original
# workbook view where TC08_VerifyProductListPLP is tabSelected="1").
$ws1.Activate()
$ws1.Range("A3:XFD7").Select()
